$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily-scrape row was inserted right before the existing row 600,
# pushing rows 600:641 down to 601:642 (dimension grows to D642).
$ws.Rows.Item(600).Insert()

# Force column A to be stored as text so the date-like string "2026/01/09"
# isn't auto-coerced into a date serial by Excel's input parser (matches
# the existing inline/shared-string date cells in this column), then reset
# the cell style back to Normal so no stray per-cell style survives.
$ws.Cells.Item(600, 1).NumberFormat = "@"
$ws.Cells.Item(600, 1).Value = "2026/01/09"
$ws.Cells.Item(600, 1).Style = "Normal"

$ws.Cells.Item(600, 2).Value = "金"
$ws.Cells.Item(600, 3).Value = 20
$ws.Cells.Item(600, 4).Value = 201
